# Merge the split "<id>...</id>" runs (tc/tcn/tl identifiers) back into a
# single run per occurrence, matching the newly downloaded tc/tcn/tl source.
#
# Before (3 runs per id):
#   <w:r><w:rPr>...Courier New, 7f6000, sz18, szCs18...</w:rPr><w:t>&lt;id&gt;</w:t></w:r>
#   <w:r><w:rPr>...</w:rPr><w:t>VALUE</w:t></w:r>
#   <w:r><w:rPr>...Courier New, 7f6000, sz18, szCs18...</w:rPr><w:t>&lt;/id&gt;</w:t></w:r>
#
# After (1 run):
#   <w:r><w:rPr>...Courier New, 7f6000, sz18, szCs18...</w:rPr><w:t>&lt;id&gt;VALUE&lt;/id&gt;</w:t></w:r>

$d = $word.ActiveDocument

function Merge-IdRun($idValue) {
    $searchText = "<id>" + $idValue + "</id>"
    $escapedText = "&lt;id&gt;" + $idValue + "&lt;/id&gt;"

    $found = $d.Content
    $found.Find.ClearFormatting()
    $found.Find.Execute($searchText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

    $start = $found.Start
    $found.Delete()

    $xml = '<?xml version="1.0" standalone="yes"?>' +
           '<?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
             '<pkg:part pkg:name="/_rels/.rels" pkg:contentType="application/vnd.openxmlformats-package.relationships+xml" pkg:padding="512">' +
               '<pkg:xmlData>' +
                 '<Relationships xmlns="http://schemas.openxmlformats.org/package/2006/relationships">' +
                   '<Relationship Id="rId1" Type="http://schemas.openxmlformats.org/officeDocument/2006/relationships/officeDocument" Target="word/document.xml"/>' +
                 '</Relationships>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
             '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
               '<pkg:xmlData>' +
                 '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
                   '<w:body>' +
                     '<w:p>' +
                       '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000">' +
                         '<w:rPr>' +
                           '<w:rFonts w:ascii="Courier New" w:eastAsia="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/>' +
                           '<w:color w:val="7f6000"/>' +
                           '<w:sz w:val="18"/>' +
                           '<w:szCs w:val="18"/>' +
                           '<w:rtl w:val="0"/>' +
                         '</w:rPr>' +
                         '<w:t xml:space="preserve">' + $escapedText + '</w:t>' +
                       '</w:r>' +
                     '</w:p>' +
                   '</w:body>' +
                 '</w:document>' +
               '</pkg:xmlData>' +
             '</pkg:part>' +
           '</pkg:package>'

    $ins = $d.Range($start, $start)
    $ins.InsertXML($xml)
}

Merge-IdRun "p058v_5"
Merge-IdRun "p059r_1"
Merge-IdRun "p059r_2"
